$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (Sheet1 -> Kedar Jadhav)
$ws.Name = "Kedar Jadhav"

# Header row. A new leading "matchNo" column is inserted, shifting every
# other header one column to the right (teamName: A->B, batterName: B->C, ...).
$headers = @("matchNo","teamName","batterName","states","runs","balls","fours","sixes","sr","opponentTeamName","venue","date","result")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Full set of match rows: the original single row is reworked (now row 5,
# matching the "23rd" match vs Chennai Super Kings) and four more rows are
# added for a total of 5 data rows.
$data = @(
    @("28th","Sunrisers Hyderabad","Kedar Jadhav","b Morris","19","19","0","1","100.00","Rajasthan Royals","Delhi","May 02","Royals won by 55 runs"),
    @("37th","Sunrisers Hyderabad","Kedar Jadhav","b Ravi Bishnoi","12","12","0","0","100.00","Punjab Kings","Sharjah","September 25","Punjab Kings won by 5 runs"),
    @("33rd","Sunrisers Hyderabad","Kedar Jadhav","lbw b Nortje","3","8","0","0","37.50","Delhi Capitals","Dubai (DSC)","September 22","Capitals won by 8 wickets (with 13 balls remaining)"),
    @("23rd","Sunrisers Hyderabad","Kedar Jadhav","","12","4","1","1","300.00","Chennai Super Kings","Delhi","April 28","Super Kings won by 7 wickets (with 9 balls remaining)"),
    @("20th","Sunrisers Hyderabad","Kedar Jadhav","st †Pant b Mishra","9","9","1","0","100.00","Delhi Capitals","Chennai","April 25","Match tied (Capitals won the one-over eliminator)")
)

# Columns E-I (runs, balls, fours, sixes, sr) hold purely numeric text
# ("19", "0", "100.00", ...) which Excel would otherwise auto-coerce to
# numbers; force them to stay text (as in the source data) with the
# standard leading-apostrophe text prefix.
$numericLooking = @(5, 6, 7, 8, 9)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Length; $c++) {
        $value = $row[$c]
        $col = $c + 1
        $cell = $ws.Cells.Item($excelRow, $col)
        if ($value -eq "") {
            $cell.Value = ""
        } elseif ($numericLooking -contains $col) {
            $cell.Value = "'" + $value
        } else {
            $cell.Value = $value
        }
    }
}
